$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 29064.13089906621
$ws.Range("D3").Value = 1258.796563159701

# Row 4 (std)
$ws.Range("B4").Value = 12371.08687036967
$ws.Range("D4").Value = 850.4644761958602

# Row 5 (min)
$ws.Range("B5").Value = 5840.039000000001

# Row 6 (25%)
$ws.Range("B6").Value = 19580.036
$ws.Range("D6").Value = 228.003

# Row 7 (50%)
$ws.Range("B7").Value = 25510.03450000008
$ws.Range("D7").Value = 1680

# Row 8 (75%)
$ws.Range("B8").Value = 40111.52275000008
$ws.Range("D8").Value = 1925

# Row 9 (max)
$ws.Range("B9").Value = 61630.61799999995
$ws.Range("D9").Value = 6230.005

# Row 10 (Total)
$ws.Range("F10").Value = 15276107200.54601

# Row 11 (Residential)
$ws.Range("G11").Value = 0.8208412111448016

# Row 12 (Community)
$ws.Range("F12").Value = 661623473.597
$ws.Range("G12").Value = 0.04331099964874242

# Row 13 (IGA)
$ws.Range("G13").Value = 0.135847789206456
